# Update the auto-date placeholders across the deck's masters/layouts.
#
# The deck's "today" stamp moved from 11/22/2022 -> 12/4/2022, which touches:
#   - the Handout Master's date placeholder   (short form, e.g. "12/4/2022")
#   - the Notes Master's date placeholder     (short form, e.g. "12/4/2022")
#   - the Slide Master's date placeholder     (long form,  e.g. "Sunday, December 4, 2022")
#   - every Slide Layout's date placeholder   (long form,  e.g. "Sunday, December 4, 2022")

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16

$shortDate = "12/4/2022"
$longDate  = "Sunday, December 4, 2022"

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shape.HasTextFrame) {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}

# Handout master (short date form, datetimeFigureOut field)
Update-DatePlaceholder $p.HandoutMaster.Shapes $shortDate

# Notes master (short date form, datetimeFigureOut field)
Update-DatePlaceholder $p.NotesMaster.Shapes $shortDate

# Slide master (long date form, datetime2 field)
Update-DatePlaceholder $p.SlideMaster.Shapes $longDate

# Every slide layout (long date form, datetime2 field)
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $longDate
}
